$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45533
# (2024-08-29) to 45534 (2024-08-30) for every data row (rows 2-28).
for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45533) {
        $cell.Value = 45534
    }
}
